# Fruta / hortaliza, semanal
# Insert two new weekly observation rows at the top of the data block
# (rows 24-25), pushing the existing rows 24-62 down to 26-64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 24.
$ws.Rows("24:25").Insert()

# --- New row 24 ---
$ws.Range("A24").Value = 5
$ws.Range("B24").Value = "Macroferia Regional de Talca"
$ws.Range("C24").Value = "Maule"
$ws.Range("D24").Value = 45079
$ws.Range("E24").Value = 7
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100107
$ws.Range("H24").Value = "Otros"
$ws.Range("I24").Value = 100107001
$ws.Range("J24").Value = "Caqui"
$ws.Range("K24").Value = "Mankaki"
$ws.Range("L24").Value = "Especial"
$ws.Range("M24").Value = 230
$ws.Range("N24").Value = 12000
$ws.Range("O24").Value = 12000
$ws.Range("P24").Value = 12000
$ws.Range("Q24").Value = "`$/caja 12 kilos granel"
$ws.Range("R24").Value = "Región de O'Higgins"
$ws.Range("S24").Value = 12000
$ws.Range("T24").Value = 1

# --- New row 25 ---
$ws.Range("A25").Value = 5
$ws.Range("B25").Value = "Macroferia Regional de Talca"
$ws.Range("C25").Value = "Maule"
$ws.Range("D25").Value = 45079
$ws.Range("E25").Value = 7
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100107
$ws.Range("H25").Value = "Otros"
$ws.Range("I25").Value = 100107001
$ws.Range("J25").Value = "Caqui"
$ws.Range("K25").Value = "Mankaki"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 210
$ws.Range("N25").Value = 10000
$ws.Range("O25").Value = 10000
$ws.Range("P25").Value = 10000
$ws.Range("Q25").Value = "`$/caja 12 kilos granel"
$ws.Range("R25").Value = "Región de O'Higgins"
$ws.Range("S25").Value = 10000
$ws.Range("T25").Value = 1
